# Update cfb_weather.xlsx with Timestamp 2025-12-28T16:23:05.311108
#
# This script mirrors a scheduled re-run of the weather-pull that feeds the
# "FBS" and "Other" sheets: the wind-direction / temperature / line-movement
# columns get refreshed numbers for a handful of games, and the Timestamp
# column (which every row shares) is bumped to the new run time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "FBS"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FBS")

# Row 2 (Georgia Southern @ Appalachian State)
$ws.Range("M2").Value = "ESE"
$ws.Range("O2").Value = 23.24
$ws.Range("P2").Value = 25.4
$ws.Range("Q2").Value = "ESE"
$ws.Range("S2").Value = -10.84
$ws.Range("T2").Value = -1.1
$ws.Range("U2").Value = 19.8

# Row 3 (Tennessee @ Illinois)
$ws.Range("O3").Value = 28.46
$ws.Range("P3").Value = 12.8
$ws.Range("S3").Value = -2.19
$ws.Range("U3").Value = 1.1

# Row 4 (Coastal Carolina @ Louisiana Tech)
$ws.Range("M4").Value = "S"
$ws.Range("O4").Value = 48.14000000000001
$ws.Range("P4").Value = 4.4
$ws.Range("U4").Value = -0.2
$ws.Range("Z4").Value = -112

# Row 5 (Iowa @ Vanderbilt)
$ws.Range("O5").Value = 45.2
$ws.Range("P5").Value = 9.300000000000001
$ws.Range("U5").Value = 1.3
$ws.Range("Y5").Value = 46.5
$ws.Range("Z5").Value = -105
$ws.Range("AB5").Value = -4
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = -2

# Row 6 (Arizona State @ Duke)
$ws.Range("O6").Value = 50
$ws.Range("P6").Value = 10.8
$ws.Range("U6").Value = 5.4

# Row 7 (Michigan @ Texas)
$ws.Range("O7").Value = 66.8
$ws.Range("P7").Value = 8.4
$ws.Range("U7").Value = -1.4
$ws.Range("AB7").Value = -7.5
$ws.Range("AF7").Value = 0

# Row 14 (Rice @ Texas State)
$ws.Range("Y14").Value = 58.5
$ws.Range("AE14").Value = -0.01680672268907563

# Row 18 (Army @ Navy)
$ws.Range("Q18").Value = "N"

# Row 20 (Old Dominion @ South Florida)
$ws.Range("Q20").Value = "WSW"

# Row 23 (Washington State @ Utah State)
$ws.Range("Q23").Value = "NNE"

# Timestamp column (AK) - every game row gets the refreshed run time
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 37).Value = "2025-12-28T16:23:05.311108"
}

# ---------------------------------------------------------------------
# Sheet "Other"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Other")

# Row 4 (Illinois State vs Montana State)
$ws2.Range("S4").Value = "SE"

# Row 5 (Villanova vs Tarleton State)
$ws2.Range("S5").Value = "NNE"

# Row 6 (South Dakota vs Montana)
$ws2.Range("S6").Value = "SE"
